$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 56, shifting existing rows 56-152 down to 57-153
$ws.Rows.Item(56).Insert()

# Populate the newly inserted row 56 with the new record
$ws.Cells.Item(56, 1).Value = 5
$ws.Cells.Item(56, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(56, 3).Value = "Maule"
$ws.Cells.Item(56, 4).Value = 45210
$ws.Cells.Item(56, 5).Value = 7
$ws.Cells.Item(56, 6).Value = 100112013
$ws.Cells.Item(56, 7).Value = "Alcachofa"
$ws.Cells.Item(56, 8).Value = "Madrigal"
$ws.Cells.Item(56, 9).Value = "Primera"
$ws.Cells.Item(56, 10).Value = 500
$ws.Cells.Item(56, 11).Value = 13000
$ws.Cells.Item(56, 12).Value = 13000
$ws.Cells.Item(56, 13).Value = 13000
$ws.Cells.Item(56, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(56, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(56, 16).Value = 325
$ws.Cells.Item(56, 17).Value = 40
$ws.Cells.Item(56, 18).Value = "Hortaliza"
